$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "TypeTest" worksheet as the last (3rd) tab, after ClassListTest.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "TypeTest"

# Row 1 - section note
$ws.Range("A1").Value = "C10"

# Row 9 - extra notes above the header row
$ws.Range("L9").Value = "ISO 8601 참고"
$ws.Range("M9").Value = "1일 10초"

# Row 10 - column headers
$ws.Range("C10").Value = "Id"
$ws.Range("D10").Value = "UIntValue"
$ws.Range("E10").Value = "ShortValue"
$ws.Range("F10").Value = "LongValue"
$ws.Range("G10").Value = "FloatValue"
$ws.Range("H10").Value = "DoubleValue"
$ws.Range("I10").Value = "CharValue"
$ws.Range("J10").Value = "StringValue"
$ws.Range("K10").Value = "EnumValue"
$ws.Range("L10").Value = "DateTimeValue"
$ws.Range("M10").Value = "TimeSpanValue"

# Row 11 - MinValue test row
$ws.Range("C11").Value = -2147483648
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = -32768
$ws.Range("F11").Value = [double]"-9.2233720368547697E+18"
$ws.Range("G11").Value = [double]"-3.4028234999999999E+38"
$ws.Range("G11").NumberFormat = "0.00E+00"
$ws.Range("H11").Formula = '="-1.7976931348623157E+308"'
$ws.Range("H11").Copy()
$ws.Range("H11").PasteSpecial(-4163)
$ws.Range("I11").Value = "0x00"
$ws.Range("J11").Value = "MinValue"
$ws.Range("K11").Value = "Sunday"
$ws.Range("L11").Value = "0001-01-01T00:00:00Z"
$ws.Range("M11").Value = "-10675199.02:48:05.4775808"

# Row 12 - sample/default value row
$ws.Range("C12").Value = 1001
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 100
$ws.Range("F12").Value = 100
$ws.Range("G12").Value = 3.1415920000000002
$ws.Range("H12").Value = 3.1415926535896999
$ws.Range("I12").Value = "a"
$ws.Range("J12").Value = "ㅋㅋㅋ"
$ws.Range("K12").Value = "Monday"
$ws.Range("L12").Value = "1986-05-26T01:05:00+09:00"
$ws.Range("M12").Value = "1.00:00:10"

# Row 13 - MaxValue test row
$ws.Range("C13").Value = 2147483647
$ws.Range("D13").Value = 4294967295
$ws.Range("E13").Value = 32767
$ws.Range("F13").Value = [double]"9.2233720368547697E+18"
$ws.Range("G13").Value = [double]"3.4028234999999999E+38"
$ws.Range("G13").NumberFormat = "0.00E+00"
$ws.Range("H13").Formula = '="1.7976931348623157E+308"'
$ws.Range("H13").Copy()
$ws.Range("H13").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("I13").Value = "0xFFFF"
$ws.Range("J13").Value = "MaxValue"
$ws.Range("K13").Value = "Saturday"
$ws.Range("L13").Value = "9999-12-31T23:59:59Z"
$ws.Range("M13").Value = "10675199.02:48:05.4775807"

# Column widths (best-effort, matching the author's manual column sizing)
$ws.Columns.Item(1).ColumnWidth = 4
$ws.Columns.Item(3).ColumnWidth = 11.86
$ws.Columns.Item(4).ColumnWidth = 10.86
$ws.Columns.Item(5).ColumnWidth = 10.29
$ws.Columns.Item(6).ColumnWidth = 13.29
$ws.Columns.Item(7).ColumnWidth = 9.86
$ws.Columns.Item(8).ColumnWidth = 25.86
$ws.Columns.Item(9).ColumnWidth = 9.71
$ws.Columns.Item(10).ColumnWidth = 10.71
$ws.Columns.Item(11).ColumnWidth = 10.57
$ws.Columns.Item(12).ColumnWidth = 25.43
$ws.Columns.Item(13).ColumnWidth = 25.71

$ws.Range("E7").Select()
